$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 360, shifting existing rows 360-441 down to 361-442.
$ws.Rows("360:360").Insert()

# Populate the new row 360 with the same constant columns as the surrounding rows,
# and the new varying values from the edit.
$ws.Range("A360").Value = 3
$ws.Range("B360").Value = 'Femacal de La Calera'
$ws.Range("C360").Value = 'Coquimbo'
$ws.Range("D360").Value = 44889
$ws.Range("E360").Value = 5
$ws.Range("F360").Value = 100114013
$ws.Range("G360").Value = 'Zanahoria'
$ws.Range("H360").Value = 'Sin especificar'
$ws.Range("I360").Value = 'Primera'
$ws.Range("J360").Value = 300
$ws.Range("K360").Value = 8500
$ws.Range("L360").Value = 9000
$ws.Range("M360").Value = 8700
$ws.Range("N360").Value = '$/saco 20 kilos'
$ws.Range("O360").Value = 'Provincia de Quillota'
$ws.Range("P360").Value = 435
$ws.Range("Q360").Value = 20
$ws.Range("R360").Value = 'Hortaliza'
